# Sync local repository with Github
# Insert a new "CD_BINARY / Global cognitive decline" row into the Variables
# table (becomes the new row 17), pushing the existing rows 17-21 down to
# rows 18-22, and place the selection/view where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 17 (PATHOLOGY_TYPE), shifting the
# remaining rows down.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row with the new variable definition.
$ws.Cells.Item(17, 1).Value = "CD_BINARY"
$ws.Cells.Item(17, 2).Value = "Factor"
$ws.Cells.Item(17, 3).Value = "NO, YES"
$ws.Cells.Item(17, 4).Value = "Global cognitive decline"
$ws.Cells.Item(17, 5).Value = "-"

# Restore the view/selection state recorded in the saved workbook.
$ws.Range("C17").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
